$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row total correct count B11: 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row B12 (corr marks) and E12 (corr/total marks text)
$ws.Range("B12").Value = 115
$ws.Range("E12").Value = "115/140"
